# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Terminal Hortofrutícola Agro Chillán"
# (Espárragos) at row 3, pushing the existing rows 3-18 down to rows 4-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 3 (row 2 = most-recent existing
# record, so the new week's record becomes the new row 3; everything that was
# row 3..18 shifts to 4..19).
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with this week's data.
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 44545
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 300000000
$ws.Range("G3").Value = "Espárragos"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = 800
$ws.Range("L3").Value = 900
$ws.Range("M3").Value = 850
$ws.Range("N3").Value = "`$/kilo"
$ws.Range("O3").Value = "Provincia de Diguillín"
$ws.Range("P3").Value = 850
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = "Hortaliza"
